$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows ("line7", "line8") are inserted right after the existing
# "line6" row. That pushes the "extr1".."extr8" rows down two places
# (from sheet rows 8-15 to rows 10-17), and two of them ("extr1"/"extr2",
# now at rows 10/11) flip their in_service flag from FALSE to TRUE.
#
# Shift the existing rows 8-15 down to rows 10-17 first, working from the
# bottom up so a source row is never clobbered before it's been copied.
# Range.Copy() (rather than writing .Value2) carries the cell style along
# with the data, which keeps column A's bold/bordered "s=1" style intact
# on the rows that move - and on the two freshly created rows below.
for ($r = 15; $r -ge 8; $r--) {
    $newRow = $r + 2
    $ws.Range("A$r`:E$r").Copy($ws.Range("A$newRow`:E$newRow"))
}

# Fill in the data for the two brand new "line7"/"line8" rows.
$ws.Range("A8").Value2 = 6
$ws.Range("B8").Value2 = "line7"
$ws.Range("C8").Value2 = 14
$ws.Range("D8").Value2 = 11
$ws.Range("E8").Value2 = $true

$ws.Range("A9").Value2 = 7
$ws.Range("B9").Value2 = "line8"
$ws.Range("C9").Value2 = 16
$ws.Range("D9").Value2 = 9
$ws.Range("E9").Value2 = $false

# "extr1" and "extr2" (now rows 10 and 11) become in_service = TRUE.
$ws.Range("E10").Value2 = $true
$ws.Range("E11").Value2 = $true
